$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (columns H and/or I) ---
$ws.Range("H274").Value = 28224
$ws.Range("I274").Value = 1324

$ws.Range("H275").Value = 28545
$ws.Range("I275").Value = 1244

$ws.Range("H277").Value = 3258

$ws.Range("H278").Value = 29958
$ws.Range("I278").Value = 2086

$ws.Range("H279").Value = 44016

$ws.Range("H280").Value = 36346
$ws.Range("I280").Value = 2413

$ws.Range("H281").Value = 45454
$ws.Range("I281").Value = 3298

$ws.Range("H285").Value = 40080

$ws.Range("H287").Value = 56639
$ws.Range("I287").Value = 3860

$ws.Range("H288").Value = 53631
$ws.Range("I288").Value = 3914

$ws.Range("H289").Value = 62503
$ws.Range("I289").Value = 3596

$ws.Range("H291").Value = 14825
$ws.Range("I291").Value = 483

# --- Append new row 292 with the latest daily stats ---
$ws.Range("A292").Value = 44186
$ws.Range("B292").Value = 155218
$ws.Range("C292").Value = 110565
$ws.Range("D292").Value = 43035
$ws.Range("E292").Value = 13289
$ws.Range("F292").Value = 2663
$ws.Range("G292").Value = 1618
$ws.Range("H292").Value = 74990
$ws.Range("I292").Value = 6520
